# Apply the PopulationParameters.xlsx edit:
# - Populate the Protein/Ontogeny columns (Q2/R2) for the existing TestPopulation row
# - Add a new population row (row 3) "TestPopulation_noOnto" which mirrors row 2
#   but does not specify Protein/Ontogeny values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Fill in Protein (Q2) and Ontogeny (R2) for the existing TestPopulation row
$ws.Range("Q2").Value = "CYP3A4, CYP2D6"
$ws.Range("R2").Value = "CYP3A4, CYP2D6"

# Add a new row for a population without ontogeny information
$ws.Range("A3").Value = "TestPopulation_noOnto"
$ws.Range("B3").Value = "Human"
$ws.Range("C3").Value = "European_ICRP_2002"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0
$ws.Range("H3").Value = "kg"
$ws.Range("K3").Value = "cm"
$ws.Range("L3").Value = 22
$ws.Range("M3").Value = 41
$ws.Range("P3").Value = "kg/m²"

$ws.Range("R3").Select()

$wb.Save()
